$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows being appended are a re-ordered repeat of the existing 8 match
# rows (2-9), appended below the current data as rows 10-17.
$rows = @(
    @(" Dubai (DSC)", " October 14 2020", "Capitals won by 13 runs", "Delhi Capitals", "Rajasthan Royals", "Ajinkya Rahane ", "2", "9", "0", "0", "22.22"),
    @(" Abu Dhabi", " October 11 2020", "Mumbai won by 5 wickets (with 2 balls remaining)", "Delhi Capitals", "Mumbai Indians", "Ajinkya Rahane ", "15", "15", "3", "0", "100.00"),
    @(" Abu Dhabi", " November 02 2020", "Capitals won by 6 wickets (with 6 balls remaining)", "Delhi Capitals", "Royal Challengers Bangalore", "Ajinkya Rahane ", "60", "46", "5", "1", "130.43"),
    @(" Dubai (DSC)", " November 05 2020", "Mumbai won by 57 runs", "Delhi Capitals", "Mumbai Indians", "Ajinkya Rahane ", "0", "3", "0", "0", "0.00"),
    @(" Dubai (DSC)", " October 27 2020", "Sunrisers won by 88 runs", "Delhi Capitals", "Sunrisers Hyderabad", "Ajinkya Rahane ", "26", "19", "3", "1", "136.84"),
    @(" Abu Dhabi", " October 24 2020", "KKR won by 59 runs", "Delhi Capitals", "Kolkata Knight Riders", "Ajinkya Rahane ", "0", "1", "0", "0", "0.00"),
    @(" Sharjah", " October 17 2020", "Capitals won by 5 wickets (with 1 ball remaining)", "Delhi Capitals", "Chennai Super Kings", "Ajinkya Rahane ", "8", "10", "1", "0", "80.00"),
    @(" Dubai (DSC)", " November 10 2020", "Mumbai won by 5 wickets (with 8 balls remaining)", "Delhi Capitals", "Mumbai Indians", "Ajinkya Rahane ", "2", "4", "0", "0", "50.00")
)

$startRow = 10
$endRow = $startRow + $rows.Count - 1

# Numeric-looking columns (G:totalRuns, H:totalBalls, I:total4s, J:total6s,
# K:sr) must stay text (matches the source data's numberStoredAsText
# convention), so format as text before writing the values.
$ws.Range("G$($startRow):K$endRow").NumberFormat = "@"

$r = $startRow
foreach ($row in $rows) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
    $ws.Range("H$r").Value = $row[7]
    $ws.Range("I$r").Value = $row[8]
    $ws.Range("J$r").Value = $row[9]
    $ws.Range("K$r").Value = $row[10]
    $r++
}
